$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C5 value from "L" to "V"
$ws.Range("C5").Value = "V"

# Update the active cell selection to C6
$ws.Range("C6").Select()
